$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New defect row (row 6)
$ws.Range("A6").Value = "D4"
$ws.Range("B6").Value = "medium"
$ws.Range("C6").Value = "open"

# E6: date of first occurrence - reuse the existing date style (same as E2:E4)
$ws.Range("E2").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Value = 41759

# F6: time of test - new number format (h:mm)
$ws.Range("F6").NumberFormat = "h:mm"
$ws.Range("F6").Value = 0.54166666666666663

$ws.Range("H6").Value = "Events that already have tweets cannot be deleted."

$excel.CutCopyMode = 0
$ws.Range("H7").Select()
